# Fruta / hortaliza, semanal
# A new weekly price record (week of 2023-11-27) is added at the top of the
# data table. All existing data rows shift down by one row, and the data
# that used to occupy the last row (row 7, "Californiana(o)") now also
# ends up duplicated into a brand-new final row (row 8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row right below the header row; this pushes every
# existing data row down by one (old row 2 -> row 3, ..., old row 7 -> row 8).
$ws.Rows("2:2").Insert()

# Make sure the newly inserted row does not inherit the bold/centered
# header formatting that Excel copies down by default, so it matches the
# plain formatting used by the rest of the data rows.
$ws.Range("A2:T2").Style = "Normal"

# Populate the new row with this week's Nispero price record.
$ws.Range("A2").Value2 = 5
$ws.Range("B2").Value2 = "Macroferia Regional de Talca"
$ws.Range("C2").Value2 = "Maule"
$ws.Range("D2").Value2 = 45257
$ws.Range("E2").Value2 = 7
$ws.Range("F2").Value2 = "Fruta"
$ws.Range("G2").Value2 = 100104
$ws.Range("H2").Value2 = "Frutos de pepita"
$ws.Range("I2").Value2 = 100104004
$ws.Range("J2").Value2 = "Níspero"
$ws.Range("K2").Value2 = "Golden Nugget"
$ws.Range("L2").Value2 = "Primera"
$ws.Range("M2").Value2 = 100
$ws.Range("N2").Value2 = 20000
$ws.Range("O2").Value2 = 20000
$ws.Range("P2").Value2 = 20000
$ws.Range("Q2").Value2 = "`$/bandeja 10 kilos"
$ws.Range("R2").Value2 = "Provincia de Limarí"
$ws.Range("S2").Value2 = 2000
$ws.Range("T2").Value2 = 10

# D is a date column formatted as a full datetime ("Fecha"); keep the same
# number format used by the other rows' D cells.
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
